$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A:W shift to B:X
$ws.Columns("A:A").Insert()

# New header cell for the inserted column (row 2 is the header row)
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Data rows 4 through 19 all get the Match ID value of 3
$ws.Range("A4:A19").ClearFormats()
$ws.Range("A4:A19").Value = 3
$ws.Range("A4:A19").Font.Bold = $true

# Row 20 (hidden summary row) also gets the value; copy formatting from row 19
# since it is the last row and direct formatting calls behave inconsistently there.
$ws.Range("A20").Value = 3
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to A2 (was W14)
$ws.Range("A2").Select()
